# Applies the "Sheets via scheduled runner" update to Chocobo_Profits.xlsx.
# Updates cached price/profit figures (columns H-N) on affected Leve rows
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, and WVR sheets.

$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(100, 8).Value2 = 200000000
$ws.Cells.Item(100, 9).Value2 = 200000000
$ws.Cells.Item(100, 10).Value2 = 0
$ws.Cells.Item(100, 11).Value2 = 200000000
$ws.Cells.Item(100, 12).Value2 = 0
$ws.Cells.Item(100, 13).Value2 = -199999459
$ws.Cells.Item(100, 14).ClearContents()
$ws.Cells.Item(106, 8).Value2 = 1797.5
$ws.Cells.Item(106, 9).Value2 = 1797.5
$ws.Cells.Item(106, 11).Value2 = 1797.5
$ws.Cells.Item(106, 13).Value2 = -1166.5
$ws.Cells.Item(138, 8).Value2 = 2477.27
$ws.Cells.Item(138, 9).Value2 = 708.25
$ws.Cells.Item(138, 10).Value2 = 2814.226
$ws.Cells.Item(138, 11).Value2 = 2124.75
$ws.Cells.Item(138, 12).Value2 = 8442.678
$ws.Cells.Item(138, 13).Value2 = 3015.25
$ws.Cells.Item(138, 14).Value2 = -18722.678
$ws.Cells.Item(141, 8).Value2 = 168531.83
$ws.Cells.Item(141, 9).Value2 = 183416.55
$ws.Cells.Item(141, 11).Value2 = 550249.6499999999
$ws.Cells.Item(141, 13).Value2 = -545069.6499999999

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(63, 8).Value2 = 11546685
$ws.Cells.Item(63, 10).Value2 = 4838.25
$ws.Cells.Item(63, 12).Value2 = 4838.25
$ws.Cells.Item(63, 14).Value2 = -6210.25
$ws.Cells.Item(66, 8).Value2 = 11546685
$ws.Cells.Item(66, 10).Value2 = 4838.25
$ws.Cells.Item(66, 12).Value2 = 24191.25
$ws.Cells.Item(66, 14).Value2 = -31055.25
$ws.Cells.Item(125, 8).Value2 = 41805.625
$ws.Cells.Item(125, 10).Value2 = 41805.625
$ws.Cells.Item(125, 12).Value2 = 41805.625
$ws.Cells.Item(125, 14).Value2 = -51645.625
$ws.Cells.Item(139, 8).Value2 = 43428.965
$ws.Cells.Item(139, 10).Value2 = 43428.965
$ws.Cells.Item(139, 12).Value2 = 43428.965
$ws.Cells.Item(139, 14).Value2 = -53708.965

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(35, 8).Value2 = 32550.8
$ws.Cells.Item(35, 10).Value2 = 32550.8
$ws.Cells.Item(35, 12).Value2 = 32550.8
$ws.Cells.Item(35, 14).Value2 = -33170.8
$ws.Cells.Item(82, 8).Value2 = 28350.555
$ws.Cells.Item(82, 9).Value2 = 14499.5
$ws.Cells.Item(82, 11).Value2 = 14499.5
$ws.Cells.Item(82, 13).Value2 = -14116.5
$ws.Cells.Item(85, 8).Value2 = 28350.555
$ws.Cells.Item(85, 9).Value2 = 14499.5
$ws.Cells.Item(85, 11).Value2 = 14499.5
$ws.Cells.Item(85, 13).Value2 = -13173.5
$ws.Cells.Item(132, 8).Value2 = 51113.332
$ws.Cells.Item(132, 10).Value2 = 51113.332
$ws.Cells.Item(132, 12).Value2 = 51113.332
$ws.Cells.Item(132, 14).Value2 = -61233.332
$ws.Cells.Item(134, 8).Value2 = 1572.35
$ws.Cells.Item(134, 9).Value2 = 1110.9445
$ws.Cells.Item(134, 10).Value2 = 5725
$ws.Cells.Item(134, 11).Value2 = 3332.8335
$ws.Cells.Item(134, 12).Value2 = 17175
$ws.Cells.Item(134, 13).Value2 = -797.8335000000002
$ws.Cells.Item(134, 14).Value2 = -22245

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value2 = 27778256
$ws.Cells.Item(16, 9).Value2 = 37037508
$ws.Cells.Item(16, 10).Value2 = 500
$ws.Cells.Item(16, 11).Value2 = 37037508
$ws.Cells.Item(16, 12).Value2 = 500
$ws.Cells.Item(16, 13).Value2 = -37037221
$ws.Cells.Item(16, 14).Value2 = -1074
$ws.Cells.Item(113, 8).Value2 = 27778256
$ws.Cells.Item(113, 9).Value2 = 37037508
$ws.Cells.Item(113, 10).Value2 = 500
$ws.Cells.Item(113, 11).Value2 = 37037508
$ws.Cells.Item(113, 12).Value2 = 500
$ws.Cells.Item(113, 13).Value2 = -37035338
$ws.Cells.Item(113, 14).Value2 = -4840
$ws.Cells.Item(134, 8).Value2 = 6748.4287
$ws.Cells.Item(134, 9).Value2 = 8094.2144
$ws.Cells.Item(134, 11).Value2 = 24282.6432
$ws.Cells.Item(134, 13).Value2 = -21747.6432
$ws.Cells.Item(138, 8).Value2 = 43790
$ws.Cells.Item(138, 10).Value2 = 43790
$ws.Cells.Item(138, 12).Value2 = 43790
$ws.Cells.Item(138, 14).Value2 = -54070
$ws.Cells.Item(140, 8).Value2 = 76557.336
$ws.Cells.Item(140, 10).Value2 = 76557.336
$ws.Cells.Item(140, 12).Value2 = 76557.336
$ws.Cells.Item(140, 14).Value2 = -86917.336
$ws.Cells.Item(141, 8).Value2 = 29033.334
$ws.Cells.Item(141, 10).Value2 = 29033.334
$ws.Cells.Item(141, 12).Value2 = 29033.334
$ws.Cells.Item(141, 14).Value2 = -39393.334

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value2 = 98.05882
$ws.Cells.Item(12, 10).Value2 = 109.71429
$ws.Cells.Item(12, 12).Value2 = 329.14287
$ws.Cells.Item(12, 14).Value2 = -675.14287
$ws.Cells.Item(106, 8).Value2 = 3679.2856
$ws.Cells.Item(106, 10).Value2 = 3679.2856
$ws.Cells.Item(106, 12).Value2 = 11037.8568
$ws.Cells.Item(106, 14).Value2 = -12929.8568
$ws.Cells.Item(129, 8).Value2 = 2750.7058
$ws.Cells.Item(129, 10).Value2 = 2282.9092
$ws.Cells.Item(129, 12).Value2 = 6848.7276
$ws.Cells.Item(129, 14).Value2 = -16848.7276

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(47, 8).Value2 = 9999.5
$ws.Cells.Item(47, 10).Value2 = 9999.5
$ws.Cells.Item(47, 12).Value2 = 9999.5
$ws.Cells.Item(47, 14).Value2 = -11135.5

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value2 = 4630
$ws.Cells.Item(7, 9).Value2 = 3852.7778
$ws.Cells.Item(7, 11).Value2 = 3852.7778
$ws.Cells.Item(7, 13).Value2 = -3740.7778
$ws.Cells.Item(40, 8).Value2 = 6829.2
$ws.Cells.Item(40, 9).Value2 = 5470.2856
$ws.Cells.Item(40, 10).Value2 = 10000
$ws.Cells.Item(40, 11).Value2 = 5470.2856
$ws.Cells.Item(40, 12).Value2 = 10000
$ws.Cells.Item(40, 13).Value2 = -5334.2856
$ws.Cells.Item(40, 14).Value2 = -10272
$ws.Cells.Item(126, 8).Value2 = 4630
$ws.Cells.Item(126, 9).Value2 = 3852.7778
$ws.Cells.Item(126, 11).Value2 = 11558.3334
$ws.Cells.Item(126, 13).Value2 = -9088.3334
$ws.Cells.Item(127, 8).Value2 = 30388.75
$ws.Cells.Item(127, 10).Value2 = 30388.75
$ws.Cells.Item(127, 12).Value2 = 30388.75
$ws.Cells.Item(127, 14).Value2 = -40308.75
$ws.Cells.Item(133, 8).Value2 = 28183.5
$ws.Cells.Item(133, 10).Value2 = 28183.5
$ws.Cells.Item(133, 12).Value2 = 28183.5
$ws.Cells.Item(133, 14).Value2 = -33243.5
$ws.Cells.Item(136, 8).Value2 = 5752.722
$ws.Cells.Item(136, 9).Value2 = 2091.5
$ws.Cells.Item(136, 10).Value2 = 7583.3335
$ws.Cells.Item(136, 11).Value2 = 6274.5
$ws.Cells.Item(136, 12).Value2 = 22750.0005
$ws.Cells.Item(136, 13).Value2 = -3724.5
$ws.Cells.Item(136, 14).Value2 = -27850.0005
$ws.Cells.Item(139, 8).Value2 = 45670
$ws.Cells.Item(139, 10).Value2 = 45670
$ws.Cells.Item(139, 12).Value2 = 45670
$ws.Cells.Item(139, 14).Value2 = -55950
$ws.Cells.Item(140, 8).Value2 = 68969.375
$ws.Cells.Item(140, 10).Value2 = 68969.375
$ws.Cells.Item(140, 12).Value2 = 68969.375
$ws.Cells.Item(140, 14).Value2 = -79329.375
$ws.Cells.Item(141, 8).Value2 = 35141.844
$ws.Cells.Item(141, 10).Value2 = 35141.844
$ws.Cells.Item(141, 12).Value2 = 35141.844
$ws.Cells.Item(141, 14).Value2 = -45501.844

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(21, 8).Value2 = 18970
$ws.Cells.Item(21, 10).Value2 = 29950
$ws.Cells.Item(21, 12).Value2 = 29950
$ws.Cells.Item(21, 14).Value2 = -30420
$ws.Cells.Item(35, 8).Value2 = 18970
$ws.Cells.Item(35, 10).Value2 = 29950
$ws.Cells.Item(35, 12).Value2 = 29950
$ws.Cells.Item(35, 14).Value2 = -30530
$ws.Cells.Item(108, 8).Value2 = 36242
$ws.Cells.Item(108, 10).Value2 = 36242
$ws.Cells.Item(108, 12).Value2 = 36242
$ws.Cells.Item(108, 14).Value2 = -43922
$ws.Cells.Item(126, 8).Value2 = 1347.6364
$ws.Cells.Item(126, 9).Value2 = 1378
$ws.Cells.Item(126, 10).Value2 = 1266.6666
$ws.Cells.Item(126, 11).Value2 = 4134
$ws.Cells.Item(126, 12).Value2 = 3799.9998
$ws.Cells.Item(126, 13).Value2 = -1664
$ws.Cells.Item(126, 14).Value2 = -8739.9998
$ws.Cells.Item(139, 8).Value2 = 35794.773
$ws.Cells.Item(139, 9).Value2 = 0
$ws.Cells.Item(139, 10).Value2 = 35794.773
$ws.Cells.Item(139, 11).Value2 = 0
$ws.Cells.Item(139, 12).Value2 = 35794.773
$ws.Cells.Item(139, 13).ClearContents()
$ws.Cells.Item(139, 14).Value2 = -46074.773
$ws.Cells.Item(140, 8).Value2 = 32381
$ws.Cells.Item(140, 10).Value2 = 32381
$ws.Cells.Item(140, 12).Value2 = 32381
$ws.Cells.Item(140, 14).Value2 = -42741
$ws.Cells.Item(141, 8).Value2 = 25334.334
$ws.Cells.Item(141, 10).Value2 = 25334.334
$ws.Cells.Item(141, 12).Value2 = 25334.334
$ws.Cells.Item(141, 14).Value2 = -35694.334

